$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update all "Förändrad" (column C) values from 45190 to 45192 for rows 2-260
$lastRow = 260
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value = 45192
    }
}

# Append new row 261 with data for case "A 44803-2023"
$newRow = 261
$ws.Cells.Item($newRow, 1).Value = "A 44803-2023"
$ws.Cells.Item($newRow, 2).Value = 45190
$ws.Cells.Item($newRow, 3).Value = 45192
$ws.Cells.Item($newRow, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($newRow, 5).Value = "SÄTER"
$ws.Cells.Item($newRow, 7).Value = 1.4
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Copy styles from row 260 (B,C date style; R wrap style)
$ws.Cells.Item(260, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(260, 3).Copy()
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122)

$ws.Cells.Item(260, 18).Copy()
$ws.Cells.Item($newRow, 18).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Match explicit row height as seen in source rows (row 260 regains it; new row 261 does not)
$ws.Rows.Item(260).RowHeight = 15
